$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Row 14 Murder ---
$ws.Range("N14").Value = -75

# --- Row 15 Rape ---
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100

# --- Row 16 Robbery ---
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 123
$ws.Range("J16").Value = 109
$ws.Range("K16").Value = 12.844036697247
$ws.Range("L16").Value = 8.849557522123
$ws.Range("M16").Value = -6.106870229007
$ws.Range("N16").Value = -52.325581395348

# --- Row 17 Fel. Assault ---
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -36.666666666666
$ws.Range("I17").Value = 191
$ws.Range("J17").Value = 179
$ws.Range("K17").Value = 6.703910614525
$ws.Range("L17").Value = 25.657894736842
$ws.Range("M17").Value = 70.535714285714
$ws.Range("N17").Value = 37.410071942446

# --- Row 18 Burglary ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 160
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 82
$ws.Range("K18").Value = 7.317073170731
$ws.Range("L18").Value = 46.666666666666
$ws.Range("M18").Value = -53.926701570680
$ws.Range("N18").Value = -82.846003898635

# --- Row 19 Gr. Larceny ---
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 111.111111111111
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 78.947368421052
$ws.Range("I19").Value = 394
$ws.Range("J19").Value = 379
$ws.Range("K19").Value = 3.957783641160
$ws.Range("L19").Value = 26.688102893890
$ws.Range("M19").Value = 52.713178294573
$ws.Range("N19").Value = 63.485477178423

# --- Row 20 G.L.A. ---
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 46
$ws.Range("H20").Value = -15.217391304347
$ws.Range("I20").Value = 229
$ws.Range("J20").Value = 327
$ws.Range("K20").Value = -29.969418960244
$ws.Range("L20").Value = 45.859872611465
$ws.Range("M20").Value = 112.037037037037
$ws.Range("N20").Value = -80.173160173160

# --- Row 21 TOTAL ---
$ws.Range("C21").Value = 44
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 46.666666666666
$ws.Range("F21").Value = 156
$ws.Range("G21").Value = 137
$ws.Range("H21").Value = 13.868613138686
$ws.Range("I21").Value = 1039
$ws.Range("J21").Value = 1087
$ws.Range("K21").Value = -4.415823367065
$ws.Range("L21").Value = 28.748451053283
$ws.Range("M21").Value = 28.113440197287
$ws.Range("N21").Value = -55.388578789179

# --- Row 23 Housing ---
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 38
$ws.Range("K23").Value = -10.526315789473
$ws.Range("L23").Value = 54.545454545454
$ws.Range("M23").Value = 41.666666666666

# --- Row 24 Petit Larceny ---
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -13.888888888888
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = -12.931034482758
$ws.Range("I24").Value = 775
$ws.Range("J24").Value = 821
$ws.Range("K24").Value = -5.602923264311
$ws.Range("L24").Value = 18.865030674846
$ws.Range("M24").Value = -6.060606060606

# --- Row 25 Retail Theft ---
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 5.882352941176
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = -3.333333333333
$ws.Range("I25").Value = 436
$ws.Range("J25").Value = 470
$ws.Range("K25").Value = -7.234042553191
$ws.Range("L25").Value = 32.121212121212

# --- Row 26 Misd. Assault ---
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 57.142857142857
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 37
$ws.Range("H26").Value = 24.324324324324
$ws.Range("I26").Value = 317
$ws.Range("J26").Value = 290
$ws.Range("K26").Value = 9.310344827586
$ws.Range("L26").Value = 19.622641509434
$ws.Range("M26").Value = 28.340080971659

# --- Row 27 UCR Rape* ---
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0

# --- Row 28 Other Sex Crimes ---
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 41
$ws.Range("K28").Value = 86.363636363636
$ws.Range("L28").Value = 28.125

# --- Row 29 Shooting Vic. ---
$ws.Range("F29").Value = "'0"
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -55.555555555555
$ws.Range("N29").Value = -66.666666666666

# --- Row 30 Shooting Inc. ---
$ws.Range("F30").Value = "'0"
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("L30").Value = -42.857142857142
$ws.Range("N30").Value = -66.666666666666

# --- Row 31 Hate Crimes ---
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100
$ws.Range("J31").Value = 3

# --- Row 33 Traffic Fatalities ---
$ws.Range("I33").Value = 4
$ws.Range("K33").Value = 100
$ws.Range("L33").Value = -20
